# IQF - Scan Btn using scanner - Row highlight without modalpopup irrespective
# of pagination & manual click pointers to func as is - Fixed
#
# This reproduces (as closely as the COM surface allows) the authoring edit:
#  * "Jumbo" sheet (sheet3) gains 6 new rows (25-30) of barcode data in the
#    B/D/F/I (label) + C/E/G/J (barcode "*label*" formula) column groups,
#    mirroring the pattern already used by rows 3-22.
#  * Row 31's first cell picks up the same "label" style used elsewhere.
#  * The "Jumbo" sheet becomes the active/selected sheet & view (zoom 70,
#    scrolled down, cell M28 selected) while "Normal" loses its prior
#    active-view state and its scroll position moves up.

$wb = $excel.ActiveWorkbook

$wsNormal = $wb.Worksheets.Item("Normal")
$wsJumbo  = $wb.Worksheets.Item("Jumbo")

# ---------------------------------------------------------------------------
# 1) "Jumbo" sheet: fill in rows 25-30 with the next batch of barcodes.
# ---------------------------------------------------------------------------

# Column groups: label column -> (first barcode value, formula column)
$colGroups = @(
    @{ Label = "B"; Formula = "C"; Values = @("JB-A00100","JB-A00101","JB-A00102","JB-A00103","JB-A00104","JB-A00105") },
    @{ Label = "D"; Formula = "E"; Values = @("JB-A00110","JB-A00111","JB-A00112","JB-A00113","JB-A00114","JB-A00115") },
    @{ Label = "F"; Formula = "G"; Values = @("JB-A00116","JB-A00117","JB-A00118","JB-A00119","JB-A00120","JB-A00121") },
    @{ Label = "I"; Formula = "J"; Values = @("JB-A00122","JB-A00123","JB-A00124","JB-A00125","JB-A00126","JB-A00127") }
)

$firstRow = 25
$lastRow  = 30

# Copy the cell formatting (fonts/fill/alignment) that rows 3-22 already use
# for each label/formula column pair down onto the new rows, one contiguous
# block at a time (keeps untouched columns, e.g. H/K/L/M, from being touched).
$wsJumbo.Range("B22:G22").Copy() | Out-Null
$wsJumbo.Range("B$firstRow`:G$lastRow").PasteSpecial(-4122) | Out-Null

$wsJumbo.Range("I22:J22").Copy() | Out-Null
$wsJumbo.Range("I$firstRow`:J$lastRow").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Match the row height already used by the populated rows above (64.5pt).
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $wsJumbo.Rows.Item($r).RowHeight = $wsJumbo.Rows.Item(22).RowHeight
}

# Write the label values column-by-column (so shared-string order matches
# the authoring order: all of B, then all of D, then all of F, then all of I)
# and the "*label*" lookup formulas alongside them.
foreach ($group in $colGroups) {
    $label = $group.Label
    $formula = $group.Formula
    $values = $group.Values

    for ($i = 0; $i -lt ($lastRow - $firstRow + 1); $i++) {
        $r = $firstRow + $i
        $wsJumbo.Range("$label$r").Value = $values[$i]
    }
    for ($i = 0; $i -lt ($lastRow - $firstRow + 1); $i++) {
        $r = $firstRow + $i
        $wsJumbo.Range("$formula$r").Formula = '="*" & ' + $label + $r + ' & "*"'
    }
}

# Row 31's first (label) cell switches from the blank "row" style to the
# "label" style (the rest of row 31 stays untouched).
$wsJumbo.Range("B22").Copy() | Out-Null
$wsJumbo.Range("B31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) View/selection changes.
# ---------------------------------------------------------------------------

# "Normal" sheet: scroll position moves up (A22 -> A10); it keeps its
# existing D10 selection, and loses the "active sheet" flag once "Jumbo" is
# activated below.
$wsNormal.Activate()
$wsNormal.Range("D10").Select() | Out-Null
$normalWindow = $excel.ActiveWindow
$normalWindow.ScrollRow = 10
$normalWindow.ScrollColumn = 1

# "Jumbo" sheet becomes the active sheet/tab, zoomed to 70%, scrolled to
# row 23, with M28 selected.
$wsJumbo.Activate()
$jumboWindow = $excel.ActiveWindow
$jumboWindow.Zoom = 70
$jumboWindow.ScrollRow = 23
$jumboWindow.ScrollColumn = 1
$wsJumbo.Range("M28").Select() | Out-Null
